$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# 1. Rewrite the "Enter all parameters value ..." Heading3 paragraph into the
#    new, re-split run structure ending with the extra "otherwise ..." text.
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Enter all parameters value")
if (-not $found) { throw "Could not find target paragraph 'Enter all parameters value...'" }

$para1 = $rng.Paragraphs(1)
$para1Range = $para1.Range

$para1Xml = '<w:p ' + $wNs + '>' `
  + '<w:pPr><w:pStyle w:val="Heading3"/></w:pPr>' `
  + '<w:r><w:t>All parameters</w:t></w:r>' `
  + '<w:r><w:t xml:space="preserve"> value</w:t></w:r>' `
  + '<w:r><w:t xml:space="preserve"> will be populated from parameter.xml file if default values set to parameters</w:t></w:r>' `
  + '<w:r><w:t>, otherwise need to enter all values manually</w:t></w:r>' `
  + '<w:r><w:t>.</w:t></w:r>' `
  + '</w:p>'

[void]$para1Range.InsertXML($para1Xml)

# ---------------------------------------------------------------------------
# 2. The following empty paragraph becomes two new paragraphs: a new
#    Heading3 paragraph describing the SecurePassword.exe step, and an
#    empty paragraph carrying the "_GoBack" bookmark.
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("All parameters value will be populated")
if (-not $found2) { throw "Could not find the rewritten parameters paragraph" }

$paraAfter = $rng2.Paragraphs(1).Next()
$paraAfterRange = $paraAfter.Range

$para2Xml = '<w:p ' + $wNs + '>' `
  + '<w:pPr><w:pStyle w:val="Heading3"/></w:pPr>' `
  + '<w:r><w:t xml:space="preserve">Encode </w:t></w:r>' `
  + '<w:proofErr w:type="spellStart"/>' `
  + '<w:r><w:t>web.config</w:t></w:r>' `
  + '<w:proofErr w:type="spellEnd"/>' `
  + '<w:r><w:t xml:space="preserve"> key values using</w:t></w:r>' `
  + '<w:r><w:t xml:space="preserve"> </w:t></w:r>' `
  + '<w:r><w:t>SecurePassword.exe</w:t></w:r>' `
  + '<w:r><w:t xml:space="preserve"> application from </w:t></w:r>' `
  + '<w:proofErr w:type="spellStart"/>' `
  + '<w:r><w:t>TextSecurity</w:t></w:r>' `
  + '<w:proofErr w:type="spellEnd"/>' `
  + '<w:r><w:t xml:space="preserve"> folder.</w:t></w:r>' `
  + '</w:p>' `
  + '<w:p ' + $wNs + '>' `
  + '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' `
  + '<w:bookmarkEnd w:id="0"/>' `
  + '</w:p>'

[void]$paraAfterRange.InsertXML($para2Xml)

# ---------------------------------------------------------------------------
# 3. Drop the stray <w:lastRenderedPageBreak/> in front of
#    "It will deploy application in IIS."
# ---------------------------------------------------------------------------
$rng3 = $d.Content
$found3 = $rng3.Find.Execute("It will deploy application in IIS.")
if (-not $found3) { throw "Could not find 'It will deploy application in IIS.' paragraph" }

$para3 = $rng3.Paragraphs(1)
$para3Range = $para3.Range

$para3Xml = '<w:p ' + $wNs + '>' `
  + '<w:pPr><w:pStyle w:val="Heading3"/></w:pPr>' `
  + '<w:r><w:t>It will deploy application in IIS.</w:t></w:r>' `
  + '</w:p>'

[void]$para3Range.InsertXML($para3Xml)
